$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Yes/No column for rows 3-7 (LandingPage, DataNexus, DataNexusHome,
# DataNexusDataLake, DataNexusConnection) from "No" to "Yes"
$ws.Range("B3:B7").Value = "Yes"

# Update the view: scroll back to the top-left corner (remove topLeftCell offset)
# and move the active selection to B2:B7 with B2 as the active cell.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B2:B7").Select()
